# Apply a cyclic shift of the data in rows 3-5:
#   old row 4 -> row 3
#   old row 5 -> row 4
#   old row 3 -> row 5
# (columns A, B, D, E, F, G, H, Q, R move with the record; the
#  "Substrat-beskrivning" note in AO follows the record that had it,
#  i.e. it moves from row 3 to row 5.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the "before" values we need, from their original rows ---

# Row 3 (will end up on row 5)
$r3_A = $ws.Range("A3").Value2
$r3_B = $ws.Range("B3").Value2
$r3_D = $ws.Range("D3").Value2
$r3_E = $ws.Range("E3").Value2
$r3_F = $ws.Range("F3").Value2
$r3_G = $ws.Range("G3").Value2
$r3_H = $ws.Range("H3").Value2
$r3_Q = $ws.Range("Q3").Value2
$r3_R = $ws.Range("R3").Value2
$r3_AO = $ws.Range("AO3").Value2

# Row 4 (will end up on row 3)
$r4_A = $ws.Range("A4").Value2
$r4_B = $ws.Range("B4").Value2
$r4_D = $ws.Range("D4").Value2
$r4_E = $ws.Range("E4").Value2
$r4_F = $ws.Range("F4").Value2
$r4_G = $ws.Range("G4").Value2
$r4_H = $ws.Range("H4").Value2
$r4_Q = $ws.Range("Q4").Value2
$r4_R = $ws.Range("R4").Value2

# Row 5 (will end up on row 4)
$r5_A = $ws.Range("A5").Value2
$r5_B = $ws.Range("B5").Value2
$r5_D = $ws.Range("D5").Value2
$r5_E = $ws.Range("E5").Value2
$r5_F = $ws.Range("F5").Value2
$r5_G = $ws.Range("G5").Value2
$r5_H = $ws.Range("H5").Value2
$r5_Q = $ws.Range("Q5").Value2
$r5_R = $ws.Range("R5").Value2

# --- write the new row 3 (old row 4's data) ---
$ws.Range("A3").Value = $r4_A
$ws.Range("B3").Value = $r4_B
$ws.Range("D3").Value = $r4_D
$ws.Range("E3").Value = $r4_E
$ws.Range("F3").Value = $r4_F
$ws.Range("G3").Value = $r4_G
$ws.Range("H3").Value = $r4_H
$ws.Range("Q3").Value = $r4_Q
$ws.Range("R3").Value = $r4_R
$ws.Range("AO3").ClearContents()

# --- write the new row 4 (old row 5's data) ---
$ws.Range("A4").Value = $r5_A
$ws.Range("B4").Value = $r5_B
$ws.Range("D4").Value = $r5_D
$ws.Range("E4").Value = $r5_E
$ws.Range("F4").Value = $r5_F
$ws.Range("G4").Value = $r5_G
$ws.Range("H4").Value = $r5_H
$ws.Range("Q4").Value = $r5_Q
$ws.Range("R4").Value = $r5_R

# --- write the new row 5 (old row 3's data) ---
$ws.Range("A5").Value = $r3_A
$ws.Range("B5").Value = $r3_B
$ws.Range("D5").Value = $r3_D
$ws.Range("E5").Value = $r3_E
$ws.Range("F5").Value = $r3_F
$ws.Range("G5").Value = $r3_G
$ws.Range("H5").Value = $r3_H
$ws.Range("Q5").Value = $r3_Q
$ws.Range("R5").Value = $r3_R
$ws.Range("AO5").Value = $r3_AO
